# Adapt the column header formatting so that the "_old" / "_new" suffixes
# used for the comparison columns are replaced with the respective AHB
# format-version identifiers ("_FV2410" for the left/old block, "_FV2504"
# for the right/new block), then turn the sheet into a proper Excel Table
# (ListObject) and freeze the header row - mirroring the upstream commit
# "chore: adapt column header formatting to respective input file names".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base names shared by both the "FV2410" (left, columns A-J) and "FV2504"
# (right, columns L-U) blocks; column K holds the unchanged "diff" header.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2410"
}

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2504"
}

# Turn the populated range into an Excel Table ("Table1") with an AutoFilter
# on the header row, as seen in the new xl/tables/table1.xml part.
$tableRange = $ws.Range("A1:U52")
$listObject = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$listObject.Name = "Table1"

# Freeze the header row (split after row 1, top-left of the scrolling pane
# is A2) so the column headers stay visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
